$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet/tab per the workbook.xml sheet name change
$ws.Name = "2023_07_05 16_32"

# Row 2
$ws.Range("D2").Value = -70
$ws.Range("I2").Value = 24597
$ws.Range("J2").Value = 40347
$ws.Range("K2").Value = 24509

# Row 3
$ws.Range("D3").Value = -74
$ws.Range("I3").Value = 23773
$ws.Range("J3").Value = 38686
$ws.Range("K3").Value = 23655

# Row 4
$ws.Range("G4").Value = 26
$ws.Range("I4").Value = 24526
$ws.Range("J4").Value = 40215

# Row 5
$ws.Range("D5").Value = -68
$ws.Range("G5").Value = 24
$ws.Range("I5").Value = 24571
$ws.Range("J5").Value = 40298
$ws.Range("K5").Value = 24477

# Row 6
$ws.Range("D6").Value = -66
$ws.Range("G6").Value = 25
$ws.Range("I6").Value = 24472
$ws.Range("J6").Value = 40101
$ws.Range("K6").Value = 24372

# Row 7
$ws.Range("D7").Value = -64
$ws.Range("I7").Value = 24616
$ws.Range("J7").Value = 40386
$ws.Range("K7").Value = 24528

# Row 8
$ws.Range("D8").Value = -62
$ws.Range("G8").Value = 24
$ws.Range("I8").Value = 24540
$ws.Range("J8").Value = 40235
$ws.Range("K8").Value = 24446

# Row 9
$ws.Range("D9").Value = -74
$ws.Range("I9").Value = 24206
$ws.Range("J9").Value = 39554
$ws.Range("K9").Value = 24112

# Row 10
$ws.Range("D10").Value = -66
$ws.Range("G10").Value = 26
$ws.Range("I10").Value = 23864
$ws.Range("J10").Value = 38864
$ws.Range("K10").Value = 23758

# Row 11
$ws.Range("D11").Value = -66
$ws.Range("I11").Value = 24311
$ws.Range("J11").Value = 39772
$ws.Range("K11").Value = 24211
